try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # 1. Rename the sheet from "licenses" to "license"
    $ws.Name = "license"

    # 2. Insert a brand-new row above row 3, which pushes the existing
    #    row 3 ("yyyyyyy" / Unity UI license) down to row 4, row 4 ("end")
    #    down to row 5, and so on through row 10 -> row 11.
    $ws.Rows.Item(3).Insert()

    # 3. The freshly inserted row 3 has no formatting yet - copy the
    #    cell formatting (styles) from the row directly below (old row 3,
    #    now row 4) so A3/B3 keep the same look (label / wrapped license text).
    $ws.Range("A4:B4").Copy()
    $ws.Range("A3:B3").PasteSpecial(-4122) # xlPasteFormats
    $excel.CutCopyMode = $false

    # 4. Populate the new row 3 with the "Easy License View" entry.
    $ws.Range("A3").Value2 = "Easy License View"
    $licenseText = "The MIT License (MIT)`n`nCopyright (c) 2018 nakagawa akihiro`n`nPermission is hereby granted, free of charge, to any person obtaining a copy`nof this software and associated documentation files (the ""Software""), to deal`nin the Software without restriction, including without limitation the rights`nto use, copy, modify, merge, publish, distribute, sublicense, and/or sell`ncopies of the Software, and to permit persons to whom the Software is`nfurnished to do so, subject to the following conditions:`n`nThe above copyright notice and this permission notice shall be included in all`ncopies or substantial portions of the Software.`n`nTHE SOFTWARE IS PROVIDED ""AS IS"", WITHOUT WARRANTY OF ANY KIND, EXPRESS OR`nIMPLIED, INCLUDING BUT NOT LIMITED TO THE WARRANTIES OF MERCHANTABILITY,`nFITNESS FOR A PARTICULAR PURPOSE AND NONINFRINGEMENT. IN NO EVENT SHALL THE`nAUTHORS OR COPYRIGHT HOLDERS BE LIABLE FOR ANY CLAIM, DAMAGES OR OTHER`nLIABILITY, WHETHER IN AN ACTION OF CONTRACT, TORT OR OTHERWISE, ARISING FROM,`nOUT OF OR IN CONNECTION WITH THE SOFTWARE OR THE USE OR OTHER DEALINGS IN THE`nSOFTWARE."
    $ws.Range("B3").Value2 = $licenseText

    # Row 3 holds as much text as row 2, so it should be exactly as tall.
    $ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

    # 5. Row 4 already kept its original height (181.5) when it shifted
    #    down from the old row 3 during the insert - nothing further to do.

    # 6. Row 5 now contains what used to live in row 4 - the "end" marker.
    #    Give A5 its "end" label (B5 stays blank).
    $ws.Range("A5").Value2 = "end"

    # 7. Move the active selection back to A1.
    [void]$ws.Range("A1").Select()

} catch {
    Write-Host "ERROR:" $_.Exception.Message
    Write-Host $_.Exception.StackTrace
}
